$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 3988.2693  # H98: 3951.1482 -> 3988.2693
$ws.Cells.Item(98, 10).Value = 4362.7856  # J98: 4271 -> 4362.7856
$ws.Cells.Item(98, 12).Value = 4362.7856  # L98: 4271 -> 4362.7856
$ws.Cells.Item(98, 14).Value = -7358.7856  # N98: -7267 -> -7358.7856

$ws.Cells.Item(104, 8).Value = 758.3333  # H104: 733.3333 -> 758.3333
$ws.Cells.Item(104, 9).Value = 737.5  # I104: 733.3333 -> 737.5
$ws.Cells.Item(104, 10).Value = 800  # J104: 0 -> 800
$ws.Cells.Item(104, 11).Value = 2212.5  # K104: 2199.9999 -> 2212.5
$ws.Cells.Item(104, 12).Value = 2400  # L104: 0 -> 2400
$ws.Cells.Item(104, 13).Value = -465.5  # M104: -452.9998999999998 -> -465.5
$ws.Cells.Item(104, 14).Value = -5894  # N104: None -> -5894

$ws.Cells.Item(122, 8).Value = 3988.2693  # H122: 3951.1482 -> 3988.2693
$ws.Cells.Item(122, 10).Value = 4362.7856  # J122: 4271 -> 4362.7856
$ws.Cells.Item(122, 12).Value = 13088.3568  # L122: 12813 -> 13088.3568
$ws.Cells.Item(122, 14).Value = -17988.3568  # N122: -17713 -> -17988.3568

$ws.Cells.Item(129, 8).Value = 1716.3  # H129: 1837.1177 -> 1716.3
$ws.Cells.Item(129, 9).Value = 1402.9  # I129: 1562 -> 1402.9
$ws.Cells.Item(129, 11).Value = 4208.700000000001  # K129: 4686 -> 4208.700000000001
$ws.Cells.Item(129, 13).Value = 791.2999999999993  # M129: 314 -> 791.2999999999993

$ws.Cells.Item(132, 8).Value = 40992.72  # H132: 44549.957 -> 40992.72
$ws.Cells.Item(132, 9).Value = 40992.72  # I132: 46529.453 -> 40992.72
$ws.Cells.Item(132, 10).Value = 0  # J132: 1001 -> 0
$ws.Cells.Item(132, 11).Value = 122978.16  # K132: 139588.359 -> 122978.16
$ws.Cells.Item(132, 12).Value = 0  # L132: 3003 -> 0
$ws.Cells.Item(132, 13).ClearContents()  # M132: remove (was -137058.359)
$ws.Cells.Item(132, 14).Value = -120448.16  # N132: -8063 -> -120448.16

$ws.Cells.Item(135, 8).Value = 1360.2106  # H135: 1360.579 -> 1360.2106
$ws.Cells.Item(135, 9).Value = 980.1667  # I135: 980.75 -> 980.1667
$ws.Cells.Item(135, 11).Value = 8821.5003  # K135: 8826.75 -> 8821.5003
$ws.Cells.Item(135, 13).Value = -6286.5003  # M135: -6291.75 -> -6286.5003

$ws.Cells.Item(137, 8).Value = 7144107  # H137: 8334574.5 -> 7144107
$ws.Cells.Item(137, 10).Value = 16667968  # J137: 25001302 -> 16667968
$ws.Cells.Item(137, 12).Value = 50003904  # L137: 75003906 -> 50003904
$ws.Cells.Item(137, 14).Value = -50009004  # N137: -75009006 -> -50009004

$ws.Cells.Item(140, 8).Value = 189990  # H140: 146660 -> 189990
$ws.Cells.Item(140, 9).Value = 0  # I140: 149990 -> 0
$ws.Cells.Item(140, 10).Value = 189990  # J140: 144995 -> 189990
$ws.Cells.Item(140, 11).Value = 0  # K140: 149990 -> 0
$ws.Cells.Item(140, 12).ClearContents()  # L140: remove (was 144995)
$ws.Cells.Item(140, 13).Value = 189990  # M140: -144810 -> 189990
$ws.Cells.Item(140, 14).Value = -200350  # N140: -155355 -> -200350

$ws.Cells.Item(141, 8).Value = 1921.75  # H141: 1987.8 -> 1921.75
$ws.Cells.Item(141, 9).Value = 1482  # I141: 1513.1666 -> 1482
$ws.Cells.Item(141, 10).Value = 5000  # J141: 2699.75 -> 5000
$ws.Cells.Item(141, 11).Value = 4446  # K141: 4539.4998 -> 4446
$ws.Cells.Item(141, 12).Value = 15000  # L141: 8099.25 -> 15000
$ws.Cells.Item(141, 13).Value = 734  # M141: 640.5002000000004 -> 734
$ws.Cells.Item(141, 14).Value = -25360  # N141: -18459.25 -> -25360

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(36, 8).Value = 9831.888999999999  # H36: 9998.375 -> 9831.888999999999
$ws.Cells.Item(36, 9).Value = 9624.5  # I36: 9999.333000000001 -> 9624.5
$ws.Cells.Item(36, 11).Value = 9624.5  # K36: 9999.333000000001 -> 9624.5
$ws.Cells.Item(36, 13).Value = -9278.5  # M36: -9653.333000000001 -> -9278.5

$ws.Cells.Item(45, 8).Value = 58557.79  # H45: 65170.47 -> 58557.79
$ws.Cells.Item(45, 9).Value = 114510.89  # I45: 146556.86 -> 114510.89
$ws.Cells.Item(45, 11).Value = 114510.89  # K45: 146556.86 -> 114510.89
$ws.Cells.Item(45, 13).Value = -114133.89  # M45: -146179.86 -> -114133.89

$ws.Cells.Item(74, 8).Value = 545226.8  # H74: 588208.3 -> 545226.8
$ws.Cells.Item(74, 9).Value = 1206.0714  # I74: 1237.4615 -> 1206.0714
$ws.Cells.Item(74, 10).Value = 827311.6  # J74: 893433.2 -> 827311.6
$ws.Cells.Item(74, 11).Value = 1206.0714  # K74: 1237.4615 -> 1206.0714
$ws.Cells.Item(74, 12).Value = 827311.6  # L74: 893433.2 -> 827311.6
$ws.Cells.Item(74, 13).Value = -332.0714  # M74: -363.4614999999999 -> -332.0714
$ws.Cells.Item(74, 14).Value = -829059.6  # N74: -895181.2 -> -829059.6

$ws.Cells.Item(77, 8).Value = 545226.8  # H77: 588208.3 -> 545226.8
$ws.Cells.Item(77, 9).Value = 1206.0714  # I77: 1237.4615 -> 1206.0714
$ws.Cells.Item(77, 10).Value = 827311.6  # J77: 893433.2 -> 827311.6
$ws.Cells.Item(77, 11).Value = 6030.357  # K77: 6187.307499999999 -> 6030.357
$ws.Cells.Item(77, 12).Value = 4136558  # L77: 4467166 -> 4136558
$ws.Cells.Item(77, 13).Value = -1662.357  # M77: -1819.307499999999 -> -1662.357
$ws.Cells.Item(77, 14).Value = -4145294  # N77: -4475902 -> -4145294

$ws.Cells.Item(88, 8).Value = 1519.9445  # H88: 1573.75 -> 1519.9445
$ws.Cells.Item(88, 9).Value = 1518.6666  # I88: 1733.25 -> 1518.6666
$ws.Cells.Item(88, 11).Value = 1518.6666  # K88: 1733.25 -> 1518.6666
$ws.Cells.Item(88, 13).Value = -1112.6666  # M88: -1327.25 -> -1112.6666

$ws.Cells.Item(91, 8).Value = 1519.9445  # H91: 1573.75 -> 1519.9445
$ws.Cells.Item(91, 9).Value = 1518.6666  # I91: 1733.25 -> 1518.6666
$ws.Cells.Item(91, 11).Value = 1518.6666  # K91: 1733.25 -> 1518.6666
$ws.Cells.Item(91, 13).Value = -114.6666  # M91: -329.25 -> -114.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(25, 8).Value = 640.25  # H25: 452 -> 640.25
$ws.Cells.Item(25, 9).Value = 520.3333  # I25: 452 -> 520.3333
$ws.Cells.Item(25, 10).Value = 1000  # J25: 0 -> 1000
$ws.Cells.Item(25, 11).Value = 520.3333  # K25: 452 -> 520.3333
$ws.Cells.Item(25, 12).Value = 1000  # L25: 0 -> 1000
$ws.Cells.Item(25, 13).Value = -285.3333  # M25: -217 -> -285.3333
$ws.Cells.Item(25, 14).Value = -1470  # N25: None -> -1470

$ws.Cells.Item(64, 8).Value = 1317.0435  # H64: 997.0625 -> 1317.0435
$ws.Cells.Item(64, 10).Value = 1600.4706  # J64: 1286.9 -> 1600.4706
$ws.Cells.Item(64, 12).Value = 1600.4706  # L64: 1286.9 -> 1600.4706
$ws.Cells.Item(64, 14).Value = -2050.4706  # N64: -1736.9 -> -2050.4706

$ws.Cells.Item(67, 8).Value = 1317.0435  # H67: 997.0625 -> 1317.0435
$ws.Cells.Item(67, 10).Value = 1600.4706  # J67: 1286.9 -> 1600.4706
$ws.Cells.Item(67, 12).Value = 1600.4706  # L67: 1286.9 -> 1600.4706
$ws.Cells.Item(67, 14).Value = -3160.4706  # N67: -2846.9 -> -3160.4706

$ws.Cells.Item(86, 8).Value = 2602.0715  # H86: 2484.8125 -> 2602.0715
$ws.Cells.Item(86, 9).Value = 1518  # I86: 1487.091 -> 1518
$ws.Cells.Item(86, 10).Value = 5312.25  # J86: 4679.8 -> 5312.25
$ws.Cells.Item(86, 11).Value = 1518  # K86: 1487.091 -> 1518
$ws.Cells.Item(86, 12).Value = 5312.25  # L86: 4679.8 -> 5312.25
$ws.Cells.Item(86, 13).Value = -395  # M86: -364.0909999999999 -> -395
$ws.Cells.Item(86, 14).Value = -7558.25  # N86: -6925.8 -> -7558.25

$ws.Cells.Item(89, 8).Value = 2602.0715  # H89: 2484.8125 -> 2602.0715
$ws.Cells.Item(89, 9).Value = 1518  # I89: 1487.091 -> 1518
$ws.Cells.Item(89, 10).Value = 5312.25  # J89: 4679.8 -> 5312.25
$ws.Cells.Item(89, 11).Value = 7590  # K89: 7435.455 -> 7590
$ws.Cells.Item(89, 12).Value = 26561.25  # L89: 23399 -> 26561.25
$ws.Cells.Item(89, 13).Value = -1974  # M89: -1819.455 -> -1974
$ws.Cells.Item(89, 14).Value = -37793.25  # N89: -34631 -> -37793.25

$ws.Cells.Item(105, 8).Value = 21587.2  # H105: 21467.2 -> 21587.2
$ws.Cells.Item(105, 9).Value = 34645.332  # I105: 26209 -> 34645.332
$ws.Cells.Item(105, 10).Value = 2000  # J105: 2500 -> 2000
$ws.Cells.Item(105, 11).Value = 34645.332  # K105: 26209 -> 34645.332
$ws.Cells.Item(105, 12).Value = 2000  # L105: 2500 -> 2000
$ws.Cells.Item(105, 13).Value = -32898.332  # M105: -24462 -> -32898.332
$ws.Cells.Item(105, 14).Value = -5494  # N105: -5994 -> -5494

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 4999  # H4: 0 -> 4999
$ws.Cells.Item(4, 9).Value = 4999  # I4: 0 -> 4999
$ws.Cells.Item(4, 11).Value = 4999  # K4: 0 -> 4999
$ws.Cells.Item(4, 13).Value = -4887  # M4: None -> -4887

$ws.Cells.Item(132, 8).Value = 65386.688  # H132: 65649.19 -> 65386.688
$ws.Cells.Item(132, 9).Value = 86156.836  # I132: 93916.55 -> 86156.836
$ws.Cells.Item(132, 10).Value = 3076.25  # J132: 3461 -> 3076.25
$ws.Cells.Item(132, 11).Value = 258470.508  # K132: 281749.65 -> 258470.508
$ws.Cells.Item(132, 12).Value = 9228.75  # L132: 10383 -> 9228.75
$ws.Cells.Item(132, 13).Value = -255940.508  # M132: -279219.65 -> -255940.508
$ws.Cells.Item(132, 14).Value = -14288.75  # N132: -15443 -> -14288.75

$ws.Cells.Item(135, 8).Value = 159990  # H135: 119993.336 -> 159990
$ws.Cells.Item(135, 10).Value = 159990  # J135: 119993.336 -> 159990
$ws.Cells.Item(135, 12).Value = 159990  # L135: 119993.336 -> 159990
$ws.Cells.Item(135, 14).Value = -170130  # N135: -130133.336 -> -170130

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 4400149.5  # H4: 4000154 -> 4400149.5
$ws.Cells.Item(4, 9).Value = 6200099.5  # I4: 5166783 -> 6200099.5
$ws.Cells.Item(4, 11).Value = 18600298.5  # K4: 15500349 -> 18600298.5
$ws.Cells.Item(4, 13).Value = -18600186.5  # M4: -15500237 -> -18600186.5

$ws.Cells.Item(5, 8).Value = 889.3570999999999  # H5: 817 -> 889.3570999999999
$ws.Cells.Item(5, 9).Value = 494.5  # I5: 457.8889 -> 494.5
$ws.Cells.Item(5, 10).Value = 1415.8334  # J5: 1625 -> 1415.8334
$ws.Cells.Item(5, 11).Value = 1483.5  # K5: 1373.6667 -> 1483.5
$ws.Cells.Item(5, 12).Value = 4247.5002  # L5: 4875 -> 4247.5002
$ws.Cells.Item(5, 13).Value = -1371.5  # M5: -1261.6667 -> -1371.5
$ws.Cells.Item(5, 14).Value = -4471.5002  # N5: -5099 -> -4471.5002

$ws.Cells.Item(12, 8).Value = 716.3889  # H12: 670.55 -> 716.3889
$ws.Cells.Item(12, 10).Value = 957.6923  # J12: 864.4 -> 957.6923
$ws.Cells.Item(12, 12).Value = 2873.0769  # L12: 2593.2 -> 2873.0769
$ws.Cells.Item(12, 14).Value = -3219.0769  # N12: -2939.2 -> -3219.0769

$ws.Cells.Item(62, 8).Value = 16260  # H62: 14569.75 -> 16260
$ws.Cells.Item(62, 9).Value = 9500  # I62: 9499.5 -> 9500
$ws.Cells.Item(62, 11).Value = 28500  # K62: 28498.5 -> 28500
$ws.Cells.Item(62, 13).Value = -27814  # M62: -27812.5 -> -27814

$ws.Cells.Item(65, 8).Value = 16260  # H65: 14569.75 -> 16260
$ws.Cells.Item(65, 9).Value = 9500  # I65: 9499.5 -> 9500
$ws.Cells.Item(65, 11).Value = 85500  # K65: 85495.5 -> 85500
$ws.Cells.Item(65, 13).Value = -82068  # M65: -82063.5 -> -82068

$ws.Cells.Item(80, 8).Value = 781.6667  # H80: 796.3333 -> 781.6667
$ws.Cells.Item(80, 9).Value = 750  # I80: 0 -> 750
$ws.Cells.Item(80, 10).Value = 797.5  # J80: 796.3333 -> 797.5
$ws.Cells.Item(80, 11).Value = 2250  # K80: 0 -> 2250
$ws.Cells.Item(80, 12).Value = 2392.5  # L80: 2388.9999 -> 2392.5
$ws.Cells.Item(80, 13).Value = -1314  # M80: None -> -1314
$ws.Cells.Item(80, 14).Value = -4264.5  # N80: -4260.9999 -> -4264.5

$ws.Cells.Item(83, 8).Value = 781.6667  # H83: 796.3333 -> 781.6667
$ws.Cells.Item(83, 9).Value = 750  # I83: 0 -> 750
$ws.Cells.Item(83, 10).Value = 797.5  # J83: 796.3333 -> 797.5
$ws.Cells.Item(83, 11).Value = 6750  # K83: 0 -> 6750
$ws.Cells.Item(83, 12).Value = 7177.5  # L83: 7166.9997 -> 7177.5
$ws.Cells.Item(83, 13).Value = -2070  # M83: None -> -2070
$ws.Cells.Item(83, 14).Value = -16537.5  # N83: -16526.9997 -> -16537.5

$ws.Cells.Item(122, 8).Value = 2526417.5  # H122: 1309.6 -> 2526417.5
$ws.Cells.Item(122, 9).Value = 10101422  # I122: 617.5 -> 10101422
$ws.Cells.Item(122, 10).Value = 1416.2222  # J122: 1482.625 -> 1416.2222
$ws.Cells.Item(122, 11).Value = 90912798  # K122: 5557.5 -> 90912798
$ws.Cells.Item(122, 12).Value = 12745.9998  # L122: 13343.625 -> 12745.9998
$ws.Cells.Item(122, 13).Value = -90910348  # M122: -3107.5 -> -90910348
$ws.Cells.Item(122, 14).Value = -17645.9998  # N122: -18243.625 -> -17645.9998

$ws.Cells.Item(135, 8).Value = 889.3570999999999  # H135: 817 -> 889.3570999999999
$ws.Cells.Item(135, 9).Value = 494.5  # I135: 457.8889 -> 494.5
$ws.Cells.Item(135, 10).Value = 1415.8334  # J135: 1625 -> 1415.8334
$ws.Cells.Item(135, 11).Value = 4450.5  # K135: 4121.0001 -> 4450.5
$ws.Cells.Item(135, 12).Value = 12742.5006  # L135: 14625 -> 12742.5006
$ws.Cells.Item(135, 13).Value = -1915.5  # M135: -1586.0001 -> -1915.5
$ws.Cells.Item(135, 14).Value = -17812.5006  # N135: -19695 -> -17812.5006

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(88, 8).Value = 0  # H88: 199999.5 -> 0
$ws.Cells.Item(88, 10).Value = 0  # J88: 199999.5 -> 0
$ws.Cells.Item(88, 12).ClearContents()  # L88: remove (was 199999.5)
$ws.Cells.Item(88, 14).Value = 0  # N88: -200901.5 -> 0

$ws.Cells.Item(91, 8).Value = 0  # H91: 199999.5 -> 0
$ws.Cells.Item(91, 10).Value = 0  # J91: 199999.5 -> 0
$ws.Cells.Item(91, 12).ClearContents()  # L91: remove (was 199999.5)
$ws.Cells.Item(91, 14).Value = 0  # N91: -203119.5 -> 0

$ws.Cells.Item(102, 8).Value = 13514700  # H102: 15152773 -> 13514700
$ws.Cells.Item(102, 9).Value = 17242520  # I102: 18519706 -> 17242520
$ws.Cells.Item(102, 10).Value = 1349.625  # J102: 1575.3334 -> 1349.625
$ws.Cells.Item(102, 11).Value = 17242520  # K102: 18519706 -> 17242520
$ws.Cells.Item(102, 12).Value = 1349.625  # L102: 1575.3334 -> 1349.625
$ws.Cells.Item(102, 13).Value = -17240898  # M102: -18518084 -> -17240898
$ws.Cells.Item(102, 14).Value = -4593.625  # N102: -4819.3334 -> -4593.625

$ws.Cells.Item(123, 8).Value = 51499.6  # H123: 51083 -> 51499.6
$ws.Cells.Item(123, 10).Value = 51499.6  # J123: 51083 -> 51499.6
$ws.Cells.Item(123, 12).Value = 51499.6  # L123: 51083 -> 51499.6
$ws.Cells.Item(123, 14).Value = -56399.6  # N123: -55983 -> -56399.6

$ws.Cells.Item(132, 8).Value = 642764.1  # H132: 642795.25 -> 642764.1
$ws.Cells.Item(132, 10).Value = 917459  # J132: 917503.4399999999 -> 917459
$ws.Cells.Item(132, 12).Value = 2752377  # L132: 2752510.32 -> 2752377
$ws.Cells.Item(132, 14).Value = -2757437  # N132: -2757570.32 -> -2757437

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 8723.35  # H7: 9103.526 -> 8723.35
$ws.Cells.Item(7, 10).Value = 10121.143  # J7: 10784.308 -> 10121.143
$ws.Cells.Item(7, 12).Value = 10121.143  # L7: 10784.308 -> 10121.143
$ws.Cells.Item(7, 14).Value = -10345.143  # N7: -11008.308 -> -10345.143

$ws.Cells.Item(42, 8).Value = 20066000  # H42: 16723000 -> 20066000
$ws.Cells.Item(42, 10).Value = 60000  # J42: 34000 -> 60000
$ws.Cells.Item(42, 12).Value = 60000  # L42: 34000 -> 60000
$ws.Cells.Item(42, 14).Value = -61126  # N42: -35126 -> -61126

$ws.Cells.Item(49, 8).Value = 20066000  # H49: 16723000 -> 20066000
$ws.Cells.Item(49, 10).Value = 60000  # J49: 34000 -> 60000
$ws.Cells.Item(49, 12).Value = 60000  # L49: 34000 -> 60000
$ws.Cells.Item(49, 14).Value = -60294  # N49: -34294 -> -60294

$ws.Cells.Item(61, 8).Value = 2595.4814  # H61: 2745.25 -> 2595.4814
$ws.Cells.Item(61, 9).Value = 2233.6191  # I61: 2334.4783 -> 2233.6191
$ws.Cells.Item(61, 10).Value = 3862  # J61: 4634.8 -> 3862
$ws.Cells.Item(61, 11).Value = 2233.6191  # K61: 2334.4783 -> 2233.6191
$ws.Cells.Item(61, 12).Value = 3862  # L61: 4634.8 -> 3862
$ws.Cells.Item(61, 13).Value = -2031.6191  # M61: -2132.4783 -> -2031.6191
$ws.Cells.Item(61, 14).Value = -4266  # N61: -5038.8 -> -4266

$ws.Cells.Item(113, 8).Value = 2595.4814  # H113: 2745.25 -> 2595.4814
$ws.Cells.Item(113, 9).Value = 2233.6191  # I113: 2334.4783 -> 2233.6191
$ws.Cells.Item(113, 10).Value = 3862  # J113: 4634.8 -> 3862
$ws.Cells.Item(113, 11).Value = 2233.6191  # K113: 2334.4783 -> 2233.6191
$ws.Cells.Item(113, 12).Value = 3862  # L113: 4634.8 -> 3862
$ws.Cells.Item(113, 13).Value = -63.61909999999989  # M113: -164.4783000000002 -> -63.61909999999989
$ws.Cells.Item(113, 14).Value = -8202  # N113: -8974.799999999999 -> -8202

$ws.Cells.Item(126, 8).Value = 8723.35  # H126: 9103.526 -> 8723.35
$ws.Cells.Item(126, 10).Value = 10121.143  # J126: 10784.308 -> 10121.143
$ws.Cells.Item(126, 12).Value = 30363.429  # L126: 32352.924 -> 30363.429
$ws.Cells.Item(126, 14).Value = -35303.429  # N126: -37292.924 -> -35303.429

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 2375.318  # H122: 2457.9524 -> 2375.318
$ws.Cells.Item(122, 9).Value = 1648.7778  # I122: 1708.1177 -> 1648.7778
$ws.Cells.Item(122, 11).Value = 4946.3334  # K122: 5124.3531 -> 4946.3334
$ws.Cells.Item(122, 13).Value = -2496.3334  # M122: -2674.3531 -> -2496.3334

$ws.Cells.Item(132, 8).Value = 1606.5227  # H132: 1619.9111 -> 1606.5227
$ws.Cells.Item(132, 10).Value = 2282.6667  # J132: 2279.318 -> 2282.6667
$ws.Cells.Item(132, 12).Value = 6848.000100000001  # L132: 6837.954000000001 -> 6848.000100000001
$ws.Cells.Item(132, 14).Value = -11908.0001  # N132: -11897.954 -> -11908.0001

$ws.Cells.Item(136, 8).Value = 56974.5  # H136: 56978.668 -> 56974.5
$ws.Cells.Item(136, 10).Value = 2823.2856  # J136: 2834 -> 2823.2856
$ws.Cells.Item(136, 12).Value = 8469.856800000001  # L136: 8502 -> 8469.856800000001
$ws.Cells.Item(136, 14).Value = -13569.8568  # N136: -13602 -> -13569.8568
